$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 2461.4
$ws.Range("I58").Value = 145
$ws.Range("J58").Value = 4005.6667
$ws.Range("K58").Value = 435
$ws.Range("L58").Value = 12017.0001
$ws.Range("M58").Value = -285
$ws.Range("N58").Value = -12317.0001
$ws.Range("H132").Value = 1001944.06
$ws.Range("I132").Value = 1949.9773
$ws.Range("J132").Value = 9801892
$ws.Range("K132").Value = 5849.9319
$ws.Range("L132").Value = 29405676
$ws.Range("M132").Value = -3319.9319
$ws.Range("N132").Value = -29410736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 789.5769
$ws.Range("I2").Value = 746
$ws.Range("J2").Value = 849
$ws.Range("K2").Value = 746
$ws.Range("L2").Value = 849
$ws.Range("M2").Value = -633
$ws.Range("N2").Value = -1075
$ws.Range("H31").Value = 7802.75
$ws.Range("I31").Value = 7802.75
$ws.Range("K31").Value = 7802.75
$ws.Range("M31").Value = -7508.75
$ws.Range("H32").Value = 4168943
$ws.Range("I32").Value = 5052439
$ws.Range("J32").Value = 3889.3572
$ws.Range("K32").Value = 5052439
$ws.Range("L32").Value = 3889.3572
$ws.Range("M32").Value = -5052152
$ws.Range("N32").Value = -4463.3572
$ws.Range("H69").Value = 50459
$ws.Range("J69").Value = 50459
$ws.Range("L69").Value = 50459
$ws.Range("N69").Value = -51957
$ws.Range("H72").Value = 50459
$ws.Range("J72").Value = 50459
$ws.Range("L72").Value = 151377
$ws.Range("N72").Value = -158865
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H116").Value = 789.5769
$ws.Range("I116").Value = 746
$ws.Range("J116").Value = 849
$ws.Range("K116").Value = 746
$ws.Range("L116").Value = 849
$ws.Range("M116").Value = 1548
$ws.Range("N116").Value = -5437
$ws.Range("H132").Value = 48512.7
$ws.Range("I132").Value = 45349.348
$ws.Range("J132").Value = 52150.55
$ws.Range("K132").Value = 136048.044
$ws.Range("L132").Value = 156451.65
$ws.Range("M132").Value = -133518.044
$ws.Range("N132").Value = -161511.65

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 789.5769
$ws.Range("I3").Value = 746
$ws.Range("J3").Value = 849
$ws.Range("K3").Value = 746
$ws.Range("L3").Value = 849
$ws.Range("M3").Value = -632
$ws.Range("N3").Value = -1077
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12387.585
$ws.Range("I31").Value = 17178.324
$ws.Range("J31").Value = 3814.6843
$ws.Range("K31").Value = 17178.324
$ws.Range("L31").Value = 3814.6843
$ws.Range("M31").Value = -16883.324
$ws.Range("N31").Value = -4404.6843
$ws.Range("H34").Value = 12387.585
$ws.Range("I34").Value = 17178.324
$ws.Range("J34").Value = 3814.6843
$ws.Range("K34").Value = 17178.324
$ws.Range("L34").Value = 3814.6843
$ws.Range("M34").Value = -16976.324
$ws.Range("N34").Value = -4218.6843

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11819072
$ws.Range("I4").Value = 4286255.5
$ws.Range("J4").Value = 25001500
$ws.Range("K4").Value = 12858766.5
$ws.Range("L4").Value = 75004500
$ws.Range("M4").Value = -12858654.5
$ws.Range("N4").Value = -75004724
$ws.Range("H68").Value = 876.2353000000001
$ws.Range("I68").Value = 463.26315
$ws.Range("J68").Value = 1121.4375
$ws.Range("K68").Value = 1389.78945
$ws.Range("L68").Value = 3364.3125
$ws.Range("M68").Value = -578.78945
$ws.Range("N68").Value = -4986.3125
$ws.Range("H71").Value = 876.2353000000001
$ws.Range("I71").Value = 463.26315
$ws.Range("J71").Value = 1121.4375
$ws.Range("K71").Value = 4169.36835
$ws.Range("L71").Value = 10092.9375
$ws.Range("M71").Value = -113.3683499999997
$ws.Range("N71").Value = -18204.9375
$ws.Range("H131").Value = 919.1739
$ws.Range("I131").Value = 516.6667
$ws.Range("J131").Value = 979.55
$ws.Range("K131").Value = 1550.0001
$ws.Range("L131").Value = 2938.65
$ws.Range("M131").Value = 3489.9999
$ws.Range("N131").Value = -13018.65

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 158.11111
$ws.Range("I2").Value = 158.11111
$ws.Range("K2").Value = 158.11111
$ws.Range("M2").Value = -45.11111
$ws.Range("H122").Value = 1921.0834
$ws.Range("I122").Value = 1473.3334
$ws.Range("J122").Value = 2368.8333
$ws.Range("K122").Value = 4420.0002
$ws.Range("L122").Value = 7106.499899999999
$ws.Range("M122").Value = -1970.0002
$ws.Range("N122").Value = -12006.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5878.5
$ws.Range("I61").Value = 6669.3335
$ws.Range("J61").Value = 5404
$ws.Range("K61").Value = 6669.3335
$ws.Range("L61").Value = 5404
$ws.Range("M61").Value = -6467.3335
$ws.Range("N61").Value = -5808
$ws.Range("H113").Value = 5878.5
$ws.Range("I113").Value = 6669.3335
$ws.Range("J113").Value = 5404
$ws.Range("K113").Value = 6669.3335
$ws.Range("L113").Value = 5404
$ws.Range("M113").Value = -4499.3335
$ws.Range("N113").Value = -9744
$ws.Range("H122").Value = 3933.7778
$ws.Range("I122").Value = 3900.5
$ws.Range("J122").Value = 4200
$ws.Range("K122").Value = 11701.5
$ws.Range("L122").Value = 12600
$ws.Range("M122").Value = -9251.5
$ws.Range("N122").Value = -17500
$ws.Range("H132").Value = 33757.78
$ws.Range("I132").Value = 1926.5294
$ws.Range("K132").Value = 5779.5882
$ws.Range("M132").Value = -3249.5882
$ws.Range("H133").Value = 30095
$ws.Range("J133").Value = 30095
$ws.Range("L133").Value = 30095
$ws.Range("N133").Value = -35155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1206.6154
$ws.Range("I113").Value = 933.6
$ws.Range("J113").Value = 2116.6667
$ws.Range("K113").Value = 2800.8
$ws.Range("L113").Value = 6350.000100000001
$ws.Range("M113").Value = -630.8000000000002
$ws.Range("N113").Value = -10690.0001
$ws.Range("H132").Value = 39045.02
$ws.Range("I132").Value = 31483.182
$ws.Range("J132").Value = 51522.05
$ws.Range("K132").Value = 94449.546
$ws.Range("L132").Value = 154566.15
$ws.Range("M132").Value = -91919.546
$ws.Range("N132").Value = -159626.15
